# Updates numeric values in Sheet1 to reflect recomputed line-flow results
# for the "380 kV" case, per the commit message: "case with 380 kV done".
# Columns C, D, E, F, G, J, L, O for rows 2-25 are updated; other columns
# (A, B, H, I, K, M, N) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$colLetters = @('C','D','E','F','G','J','L','O')
$data = @(
    ,@(2, 0.4008259790170001, 0.08476179735218814, 0.1646544839845845, 3.188501816221418, 0.002534671248454599, 0.2847810779152766, 0.1424449109323263, 9.126135068038934)
    ,@(3, 0.397770269521331, 0.08500181780032356, 0.1644652704300178, 3.139123629587729, 0.002540606040274325, 0.2833444913861243, 0.1428503099764136, 8.963443070047902)
    ,@(4, 0.3960940797312844, 0.08516658884873252, 0.1644172740985432, 3.110605269914132, 0.002544444145919403, 0.2825991955898743, 0.1431443732272442, 8.86884380537532)
    ,@(5, 0.3954612393415573, 0.08523811795603287, 0.164414856852634, 3.099435058787407, 0.002546057191693981, 0.2823298260257445, 0.1432755714425902, 8.831620452771404)
    ,@(6, 0.3953591876496176, 0.08525026032848437, 0.1644154907084037, 3.097607468376481, 0.002546328000360748, 0.2822871703623875, 0.1432980436042186, 8.825519532964677)
    ,@(7, 0.3960853417833761, 0.08516753575446323, 0.1644171720962255, 3.110452799290954, 0.002544465701410878, 0.2825954237805703, 0.1431460965786293, 8.86833643305448)
    ,@(8, 0.399730804837418, 0.08484095077689169, 0.1645750895358091, 3.171101753324621, 0.002536677379832345, 0.2842573246246474, 0.142575329664119, 9.068937128403149)
    ,@(9, 0.4084711701628123, 0.08433819011508703, 0.1654261277584936, 3.30439041496254, 0.002522936840166208, 0.2886043149128952, 0.1418137995033319, 9.50456811972532)
    ,@(10, 0.4158706283182028, 0.08405226851534309, 0.1663822655070923, 3.411187870280571, 0.002513764743439303, 0.2924660023284673, 0.141471866956131, 9.850776479271531)
    ,@(11, 0.4194509569321383, 0.08394022163348325, 0.1668893041854673, 3.461725954367154, 0.002509790191112729, 0.2943689193267147, 0.1413634587050794, 10.01404308703729)
    ,@(12, 0.420837671952313, 0.0839003758941459, 0.167091686474599, 3.481146546516129, 0.002508313404818694, 0.2951106083341699, 0.1413291771737413, 10.07670476525414)
    ,@(13, 0.4205376410056942, 0.08390884257811848, 0.1670476380776584, 3.476951363661442, 0.002508630201756692, 0.2949499329623819, 0.1413362593279004, 10.06317217017795)
    ,@(14, 0.4195644224445516, 0.08393689175511554, 0.1669057462420547, 3.463318018023926, 0.002509668129024201, 0.2944295153304495, 0.1413605027006604, 10.01918150971187)
    ,@(15, 0.4189723282730995, 0.0839544089825246, 0.16682018524914, 3.455004096420197, 0.00251030756879072, 0.2941134938964751, 0.141376233927609, 9.992345025575219)
    ,@(16, 0.415640966511944, 0.08405995302907598, 0.1663505808915531, 3.407924582323432, 0.00251402845745733, 0.2923445903930144, 0.1414798995321433, 9.840223313779688)
    ,@(17, 0.413652242413491, 0.08412931105618426, 0.1660809650593471, 3.379544867720597, 0.002516361666917743, 0.2912969238348495, 0.1415555624035036, 9.748384667770551)
    ,@(18, 0.4125285430204144, 0.08417090053958987, 0.1659326735588031, 3.36340555053269, 0.00251772230234405, 0.2907080919582867, 0.1416035190523068, 9.696104696330508)
    ,@(19, 0.4121515372986835, 0.08418527364018402, 0.1658836294877126, 3.357972595402231, 0.002518186195865904, 0.2905110844965222, 0.1416205187271835, 9.678496715856681)
    ,@(20, 0.4138618580061291, 0.08412175223929452, 0.1661089638839677, 3.382546886850747, 0.00251611136531281, 0.2914070253063841, 0.1415470487715105, 9.758104789121376)
    ,@(21, 0.419849440216268, 0.08392858296105032, 0.1669471415663288, 3.467314769812646, 0.002509362498145959, 0.2945818015051742, 0.141353198152764, 10.03207989338955)
    ,@(22, 0.4239429633025793, 0.08381739382268449, 0.1675554347012351, 3.524365370009235, 0.002505116546075089, 0.2967796970004173, 0.1412659647149219, 10.21601499353631)
    ,@(23, 0.4217416392033044, 0.08387536211031232, 0.1672252377328824, 3.493764842240381, 0.002507367663496564, 0.2955953609485107, 0.1413089149029858, 10.11739722765606)
    ,@(24, 0.4137670296160252, 0.08412516423920025, 0.1660962847029985, 3.381189124191167, 0.002516224466666168, 0.2913572064316128, 0.1415508839013313, 9.753708705221186)
    ,@(25, 0.4059355002486029, 0.08445951084074288, 0.1651378580110006, 3.266783272650969, 0.002526491120654018, 0.2873113850368156, 0.1419815707233525, 9.382154833515415)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $val = $row[$i + 1]
        $ws.Range("$col$r").Value2 = $val
    }
}
